$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold, border, centered) from existing header cell H1
# onto the two new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-22
$data = @(
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(4, 4),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(1, 1),
    @(6, 7),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(5, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
